$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "SARS_CoV_2/human/USAWA_UW415/2020"
$ws.Range("B2").Value = "QIU81213"
$ws.Range("E2").Value = 0.98422090729783
$ws.Range("F2").Value = 0.85771920387305

# Row 3
$ws.Range("A3").Value = "SARS_CoV_2/human/USA/WA_UW373/2020"
$ws.Range("B3").Value = "QIS61374"
$ws.Range("E3").Value = 0.985144503466282
$ws.Range("F3").Value = 0.85783739983794

# Row 4
$ws.Range("F4").Value = 0.857489224137931

# Row 5
$ws.Range("A5").Value = "SARS_CoV_2/human/ESP/Valencia20/2020"
$ws.Range("B5").Value = "QIU78824"
$ws.Range("E5").Value = 0.992756349952963
$ws.Range("F5").Value = 0.865286923800564

# Row 6
$ws.Range("F6").Value = 0.854837271511368

# Row 7
$ws.Range("F7").Value = 0.855448546459782

# Row 8
$ws.Range("F8").Value = 0.855957902247591

# Row 9
$ws.Range("F9").Value = 0.855957902247591

# Row 10
$ws.Range("F10").Value = 0.855129348795718

# Row 11
$ws.Range("F11").Value = 0.855524837242486

# Row 12
$ws.Range("F12").Value = 0.856185797230906

# Row 13
$ws.Range("A13").Value = "SARS_CoV_2/human/USA/WA_UW199/2020"
$ws.Range("B13").Value = "QIQ49832"
$ws.Range("E13").Value = 0.984890478319177
$ws.Range("F13").Value = 0.856057219490388

# Row 14
$ws.Range("A14").Value = "SARS_CoV_2/human/USA/WA_UW225/2020"
$ws.Range("B14").Value = "QIQ50092"
$ws.Range("E14").Value = 0.984921484653818
$ws.Range("F14").Value = 0.85608493932905

# Row 15
$ws.Range("A15").Value = "SARS_CoV_2/human/USA/WA_UW301/2020"
$ws.Range("B15").Value = "QIS60546"
$ws.Range("F15").Value = 0.855524837242486

# Row 16
$ws.Range("A16").Value = "SARS_CoV_2/human/USAWA_UW395/2020"
$ws.Range("B16").Value = "QIU80973"
$ws.Range("E16").Value = 0.984572855359372
$ws.Range("F16").Value = 0.85553772070626

# Row 17
$ws.Range("A17").Value = "SARS_CoV_2/human/USAWA_UW448/2020"
$ws.Range("B17").Value = "QIU81585"
$ws.Range("E17").Value = 0.984482297333452
$ws.Range("F17").Value = 0.855435655043253

# Row 18
$ws.Range("A18").Value = "SNU01"
$ws.Range("B18").Value = "QHZ00379"
$ws.Range("E18").Value = 0.98369129311113
$ws.Range("F18").Value = 0.854469298636485

# Row 19
$ws.Range("A19").Value = "SARS_CoV_2/human/SWE/01/2020"
$ws.Range("B19").Value = "QIC53204"
$ws.Range("E19").Value = 0.983952928590532
$ws.Range("F19").Value = 0.854684853347597

# Row 20
$ws.Range("A20").Value = "SARS_CoV_2/human/USA/WA_UW257/2020"
$ws.Range("B20").Value = "QIS30295"
$ws.Range("E20").Value = 0.983773181169757
$ws.Range("F20").Value = 0.854493580599144

# Row 21
$ws.Range("A21").Value = "SARS_CoV_2/human/USA/WA_UW340/2020"
$ws.Range("B21").Value = "QIS60978"
$ws.Range("E21").Value = 0.984129814550642
$ws.Range("F21").Value = 0.854850213980028

# Row 22
$ws.Range("A22").Value = "SARS_CoV_2/human/USA/WA_UW304/2020"
$ws.Range("B22").Value = "QIS60582"
$ws.Range("E22").Value = 0.984129814550642
$ws.Range("F22").Value = 0.854850213980028

# Row 23
$ws.Range("A23").Value = "SARS_CoV_2/human/USA/UF_2/2020"
$ws.Range("B23").Value = "QIU81910"
$ws.Range("E23").Value = 0.984397289586305
$ws.Range("F23").Value = 0.855117689015691

# Row 24
$ws.Range("A24").Value = "Australia/VIC01/2020"
$ws.Range("B24").Value = "QHR84449"
$ws.Range("E24").Value = 0.984395898350423
$ws.Range("F24").Value = 0.85510477039679

# Row 25
$ws.Range("A25").Value = "SARS_CoV_2/human/USA/WA_UW336/2020"
$ws.Range("B25").Value = "QIS60930"
$ws.Range("E25").Value = 0.984574230940704
$ws.Range("F25").Value = 0.85528310298707

# Row 26
$ws.Range("A26").Value = "SARS_CoV_2/human/USA/WA_UW261/2020"
$ws.Range("B26").Value = "QIS30335"
$ws.Range("E26").Value = 0.984485064645564
$ws.Range("F26").Value = 0.85519393669193

# Row 27
$ws.Range("F27").Value = 0.856161940431603

# Row 28
$ws.Range("A28").Value = "2019_nCoV_WHU01"
$ws.Range("B28").Value = "QHO62107"
$ws.Range("E28").Value = 0.984929552345282
$ws.Range("F28").Value = 0.855626894952737

# Row 29
$ws.Range("A29").Value = "SARS_CoV_2/human/USA/WA_NH13/2020"
$ws.Range("B29").Value = "QIS60489"
$ws.Range("E29").Value = 0.984662029605849
$ws.Range("F29").Value = 0.855359372213304

# Row 30
$ws.Range("A30").Value = "SARS_CoV_2/human/USA/WA_UW363/2020"
$ws.Range("B30").Value = "QIS61254"
$ws.Range("E30").Value = 0.984572855359372
$ws.Range("F30").Value = 0.855270197966827

# Row 31
$ws.Range("F31").Value = 0.855270197966827

# Row 32
$ws.Range("F32").Value = 0.855270197966827

# Row 33
$ws.Range("A33").Value = "SARS_CoV_2/human/USA/WA_UW370/2020"
$ws.Range("B33").Value = "QIS61338"
$ws.Range("E33").Value = 0.984483681112894
$ws.Range("F33").Value = 0.855181023720349

# Row 34
$ws.Range("A34").Value = "SARS_CoV_2/human/ESP/Valencia13/2020"
$ws.Range("B34").Value = "QIU78707"
$ws.Range("E34").Value = 0.984305332619939
$ws.Range("F34").Value = 0.855002675227394

# Row 35
$ws.Range("A35").Value = "2019_nCoV/USA_CA5/2020"
$ws.Range("B35").Value = "QHW06059"
$ws.Range("F35").Value = 0.855168108445554

# Row 36
$ws.Range("A36").Value = "SARS_CoV_2/human/CHN/KMS1/2020"
$ws.Range("B36").Value = "QIO04367"
$ws.Range("E36").Value = 0.984482297333452
$ws.Range("F36").Value = 0.855168108445554

# Row 37
$ws.Range("A37").Value = "2019_nCoV/USA_CruiseA_23/2020"
$ws.Range("B37").Value = "QIJ96493"
$ws.Range("E37").Value = 0.984302533000356
$ws.Range("F37").Value = 0.854976810560114

# Row 38
$ws.Range("A38").Value = "SARS_CoV_2/human/USA/WA_UW378/2020"
$ws.Range("B38").Value = "QIS61422"
$ws.Range("E38").Value = 0.985103915796985
$ws.Range("F38").Value = 0.855766657746855

# Row 39
$ws.Range("F39").Value = 0.855549607423269

# Row 40
$ws.Range("F40").Value = 0.85492505353319

# Row 41
$ws.Range("F41").Value = 0.855536718122601

# Row 42
$ws.Range("F42").Value = 0.855536718122601

# Row 43
$ws.Range("F43").Value = 0.854912108503613

# Row 44
$ws.Range("F44").Value = 0.854819976771196

# Row 45
$ws.Range("F45").Value = 0.85538076510547

# Row 46
$ws.Range("F46").Value = 0.854755094744369

# Row 47
$ws.Range("A47").Value = "SARS_CoV_2/human/ESP/Valencia7/2020"
$ws.Range("B47").Value = "QIQ08820"
$ws.Range("E47").Value = 0.986813186813186
$ws.Range("F47").Value = 0.857142857142857

# Row 48
$ws.Range("A48").Value = "SARS_CoV_2/human/USA/WA_UW244/2020"
$ws.Range("B48").Value = "QIS30165"
$ws.Range("E48").Value = 0.98501872659176
$ws.Range("F48").Value = 0.855270197966827

$ws.Range("K8").Select()
